# Handles float input without breaking stuff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Summary block (rows 10-12): give the label cells the same bold
# "mtitleStyle" look as the other header cells and fill in the
# actual grading numbers that were computed for this student.
# ---------------------------------------------------------------
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("A12").Style = "mtitleStyle"

$ws.Range("B10").Value = 19
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 8
$ws.Range("E10").Value = 28

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("B12").Value = 76
$ws.Range("C12").Value = -1
$ws.Range("E12").Value = "75/112"

# ---------------------------------------------------------------
# Per-question answer key (rows 16-40): fill in what the student
# actually answered in column A, coloring it green (correct),
# red (incorrect) or leaving it blank/black (not attempted) to
# match the "Correct Ans" already stored in column B.
# ---------------------------------------------------------------
$ws.Range("A16").Value = "Option A"
$ws.Range("A16").Style = "correctStyle"

$ws.Range("A17").Style = "normalStyle"

$ws.Range("A18").Value = "Option B"
$ws.Range("A18").Style = "correctStyle"

$ws.Range("A19").Value = "Option C"
$ws.Range("A19").Style = "correctStyle"

$ws.Range("A20").Style = "normalStyle"

$ws.Range("A21").Value = "Option C"
$ws.Range("A21").Style = "correctStyle"

$ws.Range("A22").Value = "Option D"
$ws.Range("A22").Style = "correctStyle"

$ws.Range("A23").Value = "Option D"
$ws.Range("A23").Style = "correctStyle"

$ws.Range("A24").Value = "Option A"
$ws.Range("A24").Style = "correctStyle"

$ws.Range("A25").Value = "Option A"
$ws.Range("A25").Style = "correctStyle"

$ws.Range("A26").Style = "normalStyle"

$ws.Range("A27").Style = "normalStyle"

$ws.Range("A28").Value = "Option D"
$ws.Range("A28").Style = "correctStyle"

$ws.Range("A29").Value = "Option D"
$ws.Range("A29").Style = "correctStyle"

$ws.Range("A30").Value = "Option B"
$ws.Range("A30").Style = "correctStyle"

$ws.Range("A31").Value = "Option B"
$ws.Range("A31").Style = "incorrectStyle"

$ws.Range("A32").Value = "Option C"
$ws.Range("A32").Style = "correctStyle"

$ws.Range("A33").Value = "Option D"
$ws.Range("A33").Style = "correctStyle"

$ws.Range("A34").Style = "normalStyle"

$ws.Range("A35").Value = "Option D"
$ws.Range("A35").Style = "correctStyle"

$ws.Range("A36").Style = "normalStyle"

$ws.Range("A37").Style = "normalStyle"

$ws.Range("A38").Value = "Option A"
$ws.Range("A38").Style = "correctStyle"

$ws.Range("A39").Value = "Option D"
$ws.Range("A39").Style = "correctStyle"

$ws.Range("A40").Style = "normalStyle"

# ---------------------------------------------------------------
# Columns D/E only keep the duplicated "Student Ans/Correct Ans"
# pair for the first couple of rows now; the rest (D19:E40) and
# the whole extra G/H pair are dropped.
# ---------------------------------------------------------------
$ws.Range("D16").Value = "Option A"
$ws.Range("D16").Style = "correctStyle"

$ws.Range("D17").Value = "Option C"
$ws.Range("D17").Style = "correctStyle"

$ws.Range("D18").Value = "Option D"
$ws.Range("D18").Style = "correctStyle"

$ws.Range("D19:E40").Clear()
$ws.Range("G15:H40").Clear()

# Drop the now-empty extra answer-key columns entirely so the
# sheet dimension shrinks back down from H to E.
$ws.Range("G1:H1").EntireColumn.Delete()
